$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Overview paragraph: drop the redundant "This is an introductory
#    exercise. " sentence (the same fact now lives, more concisely, in the
#    "Time Required" line - see change 2 below).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This is an introductory exercise. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# 2. "Time Required" line: "15-20 minutes" -> "15-20 minutes, introductory"
#    split across two runs, exactly as produced by the human edit:
#      <w:r><w:rPr/><w:tab/><w:t xml:space="preserve">15-20 minutes, </w:t></w:r>
#      <w:r><w:rPr/><w:t>introductory</w:t></w:r>
# ---------------------------------------------------------------------------

# Locate the existing run's text (tab + "15-20 minutes") and remember where
# the paragraph containing it starts, before any insertion shifts offsets.
$old = $d.Content
$old.Find.Execute("`t15-20 minutes")
$oldLen = $old.Text.Length
$oldParaStart = $old.Paragraphs.Item(1).Range.Start

# Append the replacement runs (new OOXML, via InsertXML so formatting /
# <w:tab/> / empty <w:rPr/> survive verbatim) right after the existing text.
$newRunsXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr/><w:tab/><w:t xml:space="preserve">15-20 minutes, </w:t></w:r><w:r><w:rPr/><w:t>introductory</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$old.InsertXML($newRunsXml)

# Now remove the original "<tab>15-20 minutes" run, leaving only the two
# freshly-inserted runs behind.
$toDelete = $d.Range($oldParaStart, $oldParaStart + $oldLen)
$toDelete.Delete()
